# Correct the codebook: rename the variable definition from HC_distance to HC_time
$wb = $excel.ActiveWorkbook

$codebookSheet = $wb.Worksheets.Item("Codebook")

# Fix the variable name in the codebook (was "HC_distance", should be "HC_time")
$codebookSheet.Range("A5").Value = "HC_time"

# Make the Codebook sheet the active sheet/tab, matching the saved view state
$codebookSheet.Activate()
$codebookSheet.Range("A9").Select()
